# Add a new "2020" column (Q) to the sheet, mirroring the formatting of the
# existing "2019" column (P) and filling in the new year's data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column P (rows 3-34: the thin separator row through the last data
# row) into column Q so the new column inherits the same cell styles
# (borders, fonts, number formats, alignment) as the 2019 column.
$ws.Range("P3:P34").Copy($ws.Range("Q3:Q34"))

# Header row: new year value.
$ws.Range("Q4").Value = 2020

# New 2020 data values per row (rows without an explicit value here keep the
# "-" placeholder that was copied over from column P).
$ws.Range("Q5").Value = 51
$ws.Range("Q6").Value = 29
$ws.Range("Q7").Value = 22
$ws.Range("Q8").Value = 5
$ws.Range("Q9").Value = 3
$ws.Range("Q10").Value = 2
$ws.Range("Q11").Value = 15
$ws.Range("Q12").Value = 9
$ws.Range("Q13").Value = 5
$ws.Range("Q20").Value = 7
$ws.Range("Q21").Value = 7
$ws.Range("Q26").Value = 24
$ws.Range("Q27").Value = 10
$ws.Range("Q28").Value = 14

# Rows 14-16 have real 2019 numbers in column P, but no data is available yet
# for 2020, so the new column keeps the "-" placeholder used elsewhere for
# missing values (instead of the copied-over 2019 number).
$ws.Range("Q14").Value = "-"
$ws.Range("Q15").Value = "-"
$ws.Range("Q16").Value = "-"

# Match the active selection recorded after the edit.
$ws.Range("H26").Select()
